# Regenerate the handoff report: the CI tool re-ran and produced fresh
# "Latest Handoff" timestamps (and bumped the "Ready for handoff" rows'
# Priority from "low" to "ht") for the items that were just handed off.

$wb = $excel.ActiveWorkbook

$newHandoffDatetime = "2016-08-22 22:31:17"
$newGenerateDate    = "2016-08-22 22:31:23"
$newPriority        = "ht"

# zh-cn (sheet 2) and de-de (sheet 3): rows 4-7 are the "Ready for handoff"
# rows. Column E = Priority, Column H = Latest Handoff Datetime.
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in 4..7) {
        $ws.Range("E$row").Value = $newPriority
        $ws.Range("H$row").Value = $newHandoffDatetime
    }
}

# Overview (sheet 1): rows 4-7, Column G = Latest HO Xliff Generate Date.
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in 4..7) {
    $overview.Range("G$row").Value = $newGenerateDate
}
